# Update database and change read_price algorithm
# - Remove the oldest quarter column (old column D, "Q2 1399/06") by deleting
#   the entire column, which shifts all later quarters one column to the left.
# - Append the newest quarter ("Q4 1401/12") of data in the now-empty last
#   column (M), copying the number/column formatting from the previous last
#   column (L).
# - Apply corrected figures (new "read_price" algorithm) for the quarter that
#   now lives in column I (previously column J, "Q2 1401/06").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the oldest quarter column entirely (shifts E:M -> D:L).
$ws.Range("D1:D28").EntireColumn.Delete()

# 2. Bring the new column M into existence with the same look as column L,
#    then give it the appropriate (slightly wider, "Q4") column width.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").EntireColumn.ColumnWidth = 30.17

# 3. Fill in the header / publish-date / data for the new quarter column (M).
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-12 (2)"

$ws.Range("M11").Value = 2435459
$ws.Range("M12").Value = -2093403
$ws.Range("M13").Value = 342056
$ws.Range("M14").Value = -248751
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 149968
$ws.Range("M17").Value = 243273
$ws.Range("M18").Value = -7373
$ws.Range("M19").Value = 108689
$ws.Range("M20").Value = 344589
$ws.Range("M21").Value = 215141
$ws.Range("M22").Value = 559730
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 559730
$ws.Range("M25").Value = 509
$ws.Range("M26").Value = 1100000
$ws.Range("M27").Value = 509

# 4. Corrected figures for the quarter that now sits in column I
#    (previously column J) due to the updated read_price algorithm.
$ws.Range("I9").Value = "1402-02-12 (10)"

$ws.Range("I11").Value = 1656949
$ws.Range("I12").Value = -1432696
$ws.Range("I13").Value = 224253
$ws.Range("I14").Value = -163923
$ws.Range("I16").Value = -336205
$ws.Range("I17").Value = -275875
$ws.Range("I20").Value = -201624
$ws.Range("I22").Value = -7969
$ws.Range("I24").Value = -7969
$ws.Range("I25").Value = -7
$ws.Range("I27").Value = -7
